$d = $word.ActiveDocument

# 1. Title: split "Offspring" out into a new run named "Mice"
$d.Content.Find.Execute("Adult Male Offspring", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Adult Male Mice", 2)

# 2. Affiliation fix: remove "University of " before "Michigan Medicine, Department of"
$d.Content.Find.Execute("University of Michigan Medicine, Department of ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Michigan Medicine, Department of ", 2)

# 3. Affiliation fix: remove stray comma "Division of, Diabetes" -> "Division of Diabetes"
$d.Content.Find.Execute("Division of, Diabetes", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Division of Diabetes", 2)
